$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Swap rows 3 and 4 (Paula Andrea Martinez Castro <-> Juan Carlos Diaz Perez),
#    including their "Fecha expiracion Licencia" (column F) values which travel with the row.
$ws.Range("A3").Value = "Juan"
$ws.Range("B3").Value = "Carlos"
$ws.Range("C3").Value = "Diaz Perez"
$ws.Range("D3").Value = 105369876
$ws.Range("E3").Value = 555553
$ws.Range("F3").Value = "'2019-03-22"

$ws.Range("A4").Value = "Paula"
$ws.Range("B4").Value = "Andrea"
$ws.Range("C4").Value = "Martinez Castro"
$ws.Range("D4").Value = 105369875
$ws.Range("E4").Value = 555552
$ws.Range("F4").Value = "'2019-02-22"

# 2. Insert a new column at G. Because column G currently shares its custom width
#    with column F, the new blank column inherits that exact custom width, and the
#    "col" range in the worksheet XML extends to cover the new column as well.
$ws.Columns("G").Insert()

# At this point: G is blank (new), H holds what used to be "Fecha de Nacimiento",
# I holds what used to be "Estado". Move the "Fecha de Nacimiento" column back to G
# so the new blank (future "Genero") column ends up at H, matching the target layout.
$ws.Range("H1:H7").Cut()
$ws.Range("G1").PasteSpecial(-4104)

# 3. Mark Laura's record (row 2) as used.
$ws.Range("I2").Value = "Usado"

# 4. Populate the new "Genero" column (H).
$ws.Range("H1").Value = "Genero"
$ws.Range("H2").Value = "'Mujer"
$ws.Range("H3").Value = "'Hombre"
$ws.Range("H4").Value = "'Mujer"
$ws.Range("H5").Value = "'Hombre"
$ws.Range("H6").Value = "'Hombre"
$ws.Range("H7").Value = "'Hombre"

# 5. Update the remembered selection to match the target workbook view.
$ws.Range("G20").Select()
